# Update cryptocurrency price (D) and 1h volume change (E) columns
# with freshly scraped figures from the GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells are plain numeric-looking text (e.g. "214.97", "0.0622").
# Format them as Text first so Excel keeps the exact scraped string
# instead of silently converting it to a Number.
$priceCells = @("D2", "D3", "D5", "D10", "D11", "D12", "D13", "D15", "D16", "D17", "D18", "D19", "D25", "D29", "D32", "D36", "D38", "D40", "D41", "D45", "D46", "D47", "D49", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.121.47"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").Value = "1.678.80"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "214.97"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +1.74%  "
$ws.Range("E9").Value = "  +5.24%  "
$ws.Range("D10").Value = "0.0622"
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("D11").Value = "0.0884"
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("D12").Value = "1.916.98"
$ws.Range("E12").Value = "  +0.30%  "
$ws.Range("D13").Value = "1.694.66"
$ws.Range("E13").Value = "  +1.09%  "
$ws.Range("E14").Value = "  +1.42%  "
$ws.Range("D15").Value = "0.535"
$ws.Range("E15").Value = "  +1.62%  "
$ws.Range("D16").Value = "66.15"
$ws.Range("E16").Value = "  +0.74%  "
$ws.Range("D17").Value = "27.118.79"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("D18").Value = "238.37"
$ws.Range("E18").Value = "  +1.37%  "
$ws.Range("D19").Value = "8.12"
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("E22").Value = "  +1.35%  "
$ws.Range("E23").Value = "  +3.23%  "
$ws.Range("E24").Value = "  -2.86%  "
$ws.Range("D25").Value = "146.69"
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("E26").Value = "  +0.76%  "
$ws.Range("E27").Value = "  +1.35%  "
$ws.Range("E28").Value = "  +0.35%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E30").Value = "  +0.35%  "
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("D32").Value = "1.559.64"
$ws.Range("E32").Value = "  +5.06%  "
$ws.Range("E33").Value = "  +0.92%  "
$ws.Range("E34").Value = "  +1.72%  "
$ws.Range("E35").Value = "  +1.21%  "
$ws.Range("D36").Value = "0.602"
$ws.Range("E36").Value = "  +2.92%  "
$ws.Range("E37").Value = "  -1.14%  "
$ws.Range("D38").Value = "0.935"
$ws.Range("E38").Value = "  +4.15%  "
$ws.Range("E39").Value = "  +2.60%  "
$ws.Range("D40").Value = "1.05"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").Value = "68.98"
$ws.Range("E41").Value = "  +2.11%  "
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("E43").Value = "  -3.19%  "
$ws.Range("E44").Value = "  -1.66%  "
$ws.Range("D45").Value = "1.825.47"
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("D46").Value = "0.782"
$ws.Range("E46").Value = "  +0.17%  "
$ws.Range("D47").Value = "90.76"
$ws.Range("E47").Value = "  +0.17%  "
$ws.Range("E48").Value = "  +2.95%  "
$ws.Range("D49").Value = "0.0₆0106"
$ws.Range("E49").Value = "  +0.52%  "
$ws.Range("D50").Value = "0.104"
$ws.Range("E50").Value = "  +2.55%  "
$ws.Range("D51").Value = "8.12"
$ws.Range("E51").Value = "  +4.87%  "
